$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(65).Delete()

$ws.Range("B3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B3"), "https://http2.mlstatic.com/D_NQ_NP_2X_744049-MLA99539379058_122025-F.webp")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://mercadolivre.com/sec/1Fw8u2p")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://http2.mlstatic.com/D_Q_NP_953439-MLA99450372734_112025-F.webp")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://mercadolivre.com/sec/1HbZVwu")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://http2.mlstatic.com/D_Q_NP_801309-MLB89093855460_082025-F-xicara-copo-cafe-termico-inox-80ml-expresso-parede-dupla.webp")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://http2.mlstatic.com/D_NQ_NP_2X_709939-MLB89343355425_082025-F-kit-5-arandela-meia-lua-led-8w-8-fachos-ip66-bivolt-3000k.webp")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://http2.mlstatic.com/D_Q_NP_900674-MLA100087555981_122025-F.webp")
$ws.Hyperlinks.Add($ws.Range("A198"), "https://www.mercadolivre.com.br/importados/compra-internacional", "redirect=landing_international&origin=vip", "", "https://www.mercadolivre.com.br/importados/compra-internacional - redirect=landing_international&origin=vip")
$ws.Range("A198").Value = "Ver produtos Internacional"
$ws.Range("A198").WrapText = $true

Write-Host ("Hyperlinks count: " + $ws.Hyperlinks.Count)
